$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of rows 5 and 6 across columns A:E
$a5 = $ws.Range("A5").Value()
$b5 = $ws.Range("B5").Value()
$c5 = $ws.Range("C5").Value()
$d5 = $ws.Range("D5").Value()
$e5 = $ws.Range("E5").Value()

$a6 = $ws.Range("A6").Value()
$b6 = $ws.Range("B6").Value()
$c6 = $ws.Range("C6").Value()
$d6 = $ws.Range("D6").Value()
$e6 = $ws.Range("E6").Value()

$ws.Range("A5").Value = $a6
$ws.Range("B5").Value = $b6
$ws.Range("C5").Value = $c6
$ws.Range("D5").Value = $d6
$ws.Range("E5").Value = $e6

$ws.Range("A6").Value = $a5
$ws.Range("B6").Value = $b5
$ws.Range("C6").Value = $c5
$ws.Range("D6").Value = $d5
$ws.Range("E6").Value = $e5
